$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InfoProfessionel")

# Mapping of row -> (Statut, Poste)
$data = @{
    2  = @("pascadre", "Mecanicien")
    3  = @("pas cadre", "agent de securite")
    4  = @("cadre", "Informaticien")
    5  = @(" cadre", "developpeur")
    6  = @("cadre", "developpeur")
    7  = @("pas cadre", "Mecanicien")
    8  = @("pas cadre", "Mecanicien")
    9  = @("pas cadre", "Mecanicien")
    10 = @("pas cadre", "Mecanicien")
    11 = @("pas cadre", "Mecanicien")
    12 = @("pas cadre", "Mecanicien")
    13 = @("pas cadre", "Mecanicien")
    14 = @("pas cadre", "Mecanicien")
    15 = @("pas cadre", "Mecanicien")
    16 = @("pas cadre", "Mecanicien")
    17 = @("pas cadre", "Mecanicien")
    18 = @("pas cadre", "Mecanicien")
    19 = @("cadre", "gestionnaire comptable")
    20 = @("pas cadre", "agent d entretien")
    21 = @("cadre", "gestionnaire production")
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
}

# Adjust column widths to fit new content
$ws.Columns.Item(3).ColumnWidth = 9.28515625
$ws.Columns.Item(4).ColumnWidth = 22.28515625

# Update selection to match the diff (D7:D18, active cell D7)
$ws.Range("D7:D18").Select() | Out-Null
